$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86:109 down to 87:110
$ws.Rows(86).Insert()

# Populate the newly inserted row 86 with the new daily price record
$ws.Cells.Item(86, 1).Value = 10
$ws.Cells.Item(86, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(86, 3).Value = "La Araucanía"
$ws.Cells.Item(86, 4).Value = 45229
$ws.Cells.Item(86, 5).Value = 9
$ws.Cells.Item(86, 6).Value = 100112042
$ws.Cells.Item(86, 7).Value = "Locoto"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 60
$ws.Cells.Item(86, 11).Value = 3800
$ws.Cells.Item(86, 12).Value = 3800
$ws.Cells.Item(86, 13).Value = 3800
$ws.Cells.Item(86, 14).Value = "$/kilo"
$ws.Cells.Item(86, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(86, 16).Value = 3800
$ws.Cells.Item(86, 17).Value = 1
$ws.Cells.Item(86, 18).Value = "Hortaliza"
